$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Fill in the missing test case name for row 3 (TestCase2), matching the
# sequence TestCase1 (B2), TestCase2 (B3), TestCase3 (B4), TestCase4 (B5)
$ws.Range("B3").Value = "TestCase2"

# Update the active selection on the sheet to B6
$ws.Range("B6").Select()
